# Loan RBI, Variable Instalments
#
# The "Repayment Schedule" sheet gains a new (blank) column just before the
# old "Late" column: the old N/O/P columns ("Late", "Original"/"Heading",
# "Over Due"/"Outstanding") shift right into O/P/Q, leaving a fresh, empty
# column N behind. Insert a whole column at column N (14) so every row
# (header + 13 instalment rows) shifts in one go, preserving styles/values.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")

$ws.Columns.Item(14).Insert()

# The author was last working on the "Repayment Schedule" tab, with cell
# T7 selected, when the workbook was saved.
$ws.Activate()
$ws.Range("T7").Select()
